# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column E = Periodo Mora (text), Column F = Valor Mora, Column G = Salario Basico
# Rows 16-17: IVAN RAFAEL REDONDO FIGUEROA periods re-ordered chronologically (1803, 1804)
# Rows 18-30: YEISON ENRIQUE JULIO PARRA periods re-ordered chronologically (1808..1907)
#             and Salario Basico updated 781242 -> 828116

$rows = @(
    @{ Row = 16; Periodo = "1803"; Mora = 31249; Salario = 781242 },
    @{ Row = 17; Periodo = "1804"; Mora = 31249; Salario = 781242 },
    @{ Row = 18; Periodo = "1808"; Mora = 31249; Salario = 828116 },
    @{ Row = 19; Periodo = "1809"; Mora = 31249; Salario = 828116 },
    @{ Row = 20; Periodo = "1810"; Mora = 31249; Salario = 828116 },
    @{ Row = 21; Periodo = "1811"; Mora = 31249; Salario = 828116 },
    @{ Row = 22; Periodo = "1812"; Mora = 31249; Salario = 828116 },
    @{ Row = 23; Periodo = "1901"; Mora = 31249; Salario = 828116 },
    @{ Row = 24; Periodo = "1902"; Mora = 31249; Salario = 828116 },
    @{ Row = 25; Periodo = "1903"; Mora = 31249; Salario = 828116 },
    @{ Row = 26; Periodo = "1904"; Mora = 31249; Salario = 828116 },
    @{ Row = 27; Periodo = "1905"; Mora = 31249; Salario = 828116 },
    @{ Row = 28; Periodo = "1906"; Mora = 33125; Salario = 828116 },
    @{ Row = 29; Periodo = "1907"; Mora = 20979; Salario = 828116 },
    @{ Row = 30; Periodo = "1907"; Mora = 1104;  Salario = 828116 }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 5).Value = $r.Periodo
    $ws.Cells.Item($r.Row, 6).Value = $r.Mora
    $ws.Cells.Item($r.Row, 7).Value = $r.Salario
}
